# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) holds freshly recalculated values. Write the
# updated values for each affected row directly into column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    3  = 0
    4  = 0
    5  = 1
    6  = 3
    7  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    18 = 1
    20 = 1
    21 = 0
    22 = 1
    24 = 1
    25 = 0
    26 = 1
    27 = 1
    28 = 0
    29 = 1
    30 = 0
    31 = 3
    32 = 0
    33 = 1
    34 = 1
    35 = 2
    36 = 0
    37 = 1
    38 = 3
    39 = 1
    40 = 3
    41 = 1
    42 = 0
    43 = 2
    44 = 1
    45 = 1
    46 = 2
    47 = 0
    48 = 0
    49 = 1
    50 = 2
    51 = 0
    52 = 1
    53 = 0
    54 = 2
    55 = 2
    56 = 1
    57 = 3
    59 = 1
    60 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
